$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.726722478866577
$ws.Range("C3").Value = 0.08533525466918945
$ws.Range("C4").Value = 0.3453719615936279
$ws.Range("C5").Value = 0.05210423469543457
$ws.Range("C6").Value = 0.1052999496459961
$ws.Range("C7").Value = 1.38190770149231
$ws.Range("C8").Value = 1.689863443374634
$ws.Range("C9").Value = 1.378387928009033
$ws.Range("C10").Value = 0.3627090454101562
$ws.Range("C11").Value = 0.1477208137512207
$ws.Range("C12").Value = 0.2122530937194824
$ws.Range("C13").Value = 0.1777474880218506
$ws.Range("C14").Value = 0.2160978317260742
$ws.Range("C15").Value = 0.1105682849884033
$ws.Range("C16").Value = 0.2845001220703125
$ws.Range("C17").Value = 0.07275986671447754
$ws.Range("C18").Value = 0.5395534038543701
$ws.Range("C19").Value = 0.3739371299743652
$ws.Range("C20").Value = 0.4800035953521729
$ws.Range("C21").Value = 0.1672422885894775
$ws.Range("C22").Value = 0.3854775428771973
$ws.Range("C23").Value = 0.2192008495330811
$ws.Range("C24").Value = 0.2800445556640625
$ws.Range("C25").Value = 0.3421554565429688
$ws.Range("C26").Value = 1.173881530761719
$ws.Range("C27").Value = 0.4519391059875488
$ws.Range("C28").Value = 0.3369958400726318
$ws.Range("C29").Value = 0.8598229885101318
$ws.Range("C30").Value = 1.0510094165802
$ws.Range("C31").Value = 0.2529809474945068
$ws.Range("C32").Value = 0.7779154777526855
$ws.Range("C33").Value = 0.1767613887786865
$ws.Range("C34").Value = 0.4078166484832764
$ws.Range("C35").Value = 5.457901239395142
$ws.Range("C36").Value = 5.517557859420776
$ws.Range("C37").Value = 4.731856822967529
$ws.Range("C38").Value = 0.4133286476135254
$ws.Range("C39").Value = 0.4314947128295898
$ws.Range("C40").Value = 0.6011745929718018
$ws.Range("C41").Value = 0.5564999580383301
$ws.Range("C42").Value = 0.593944787979126
$ws.Range("C43").Value = 0.313683032989502
$ws.Range("C44").Value = 0.7046561241149902
$ws.Range("C45").Value = 0.1247367858886719
$ws.Range("C46").Value = 1.847059488296509
$ws.Range("C47").Value = 1.476282596588135
$ws.Range("C48").Value = 1.168903589248657
$ws.Range("C49").Value = 0.4819369316101074
$ws.Range("C50").Value = 1.190032243728638
$ws.Range("C51").Value = 0.5104036331176758
$ws.Range("C52").Value = 0.6044430732727051
$ws.Range("C53").Value = 0.9400427341461182
$ws.Range("C54").Value = 4.555555820465088
$ws.Range("C55").Value = 0.8267111778259277
$ws.Range("C56").Value = 0.731238842010498
$ws.Range("C57").Value = 3.479123830795288
$ws.Range("C58").Value = 6.136823415756226
$ws.Range("C59").Value = 1.567386150360107
$ws.Range("C60").Value = 5.232789754867554
$ws.Range("C61").Value = 1.177037000656128
$ws.Range("C62").Value = 2.754518032073975
$ws.Range("C63").Value = 42.53421807289124
$ws.Range("C64").Value = 43.66791987419128
$ws.Range("C65").Value = 44.12336611747742
$ws.Range("C66").Value = 2.827271938323975
$ws.Range("C67").Value = 2.953197002410889
$ws.Range("C68").Value = 4.136046648025513
$ws.Range("C69").Value = 3.130542516708374
$ws.Range("C70").Value = 4.131512880325317
$ws.Range("C71").Value = 1.933137893676758
$ws.Range("C72").Value = 5.557842254638672
$ws.Range("C73").Value = 1.162195920944214
$ws.Range("C74").Value = 12.45959091186523
$ws.Range("C75").Value = 8.528643846511841
$ws.Range("C76").Value = 8.269576549530029
$ws.Range("C77").Value = 3.101171493530273
$ws.Range("C78").Value = 7.867181539535522
$ws.Range("C79").Value = 3.781141757965088
$ws.Range("C80").Value = 5.211909532546997
$ws.Range("C81").Value = 6.173851490020752
$ws.Range("C82").Value = 42.18368887901306
$ws.Range("C83").Value = 5.18602442741394
$ws.Range("C84").Value = 5.230589389801025
$ws.Range("C85").Value = 32.00463390350342
$ws.Range("A86").Value = '''''que'''
$ws.Range("A87").Value = '''''que'''
$ws.Range("C87").Value = 0.0000007152557373046875
$ws.Range("A88").Value = '''''que'''
$ws.Range("C88").Value = 0.000000476837158203125
$ws.Range("A89").Value = '''''que'''
$ws.Range("A90").Value = '''''que'''
$ws.Range("C90").Value = 0.000000476837158203125
$ws.Range("A91").Value = '''''que'''
$ws.Range("A92").Value = '''''que'''
$ws.Range("C92").Value = 0.0000007152557373046875
$ws.Range("A93").Value = '''''que'''
$ws.Range("C93").Value = 0.0000007152557373046875
$ws.Range("A94").Value = '''''que'''
$ws.Range("A95").Value = '''''que'''
$ws.Range("C95").Value = 0.0000007152557373046875
$ws.Range("A96").Value = '''''que'''
$ws.Range("C96").Value = 0.0000007152557373046875
$ws.Range("A97").Value = '''''que'''
$ws.Range("C97").Value = 0.00000095367431640625
$ws.Range("A98").Value = '''''que'''
$ws.Range("A99").Value = '''''que'''
$ws.Range("C99").Value = 0.0000007152557373046875
$ws.Range("A100").Value = '''''que'''
$ws.Range("C100").Value = 0.00000095367431640625
$ws.Range("A101").Value = '''''que'''
$ws.Range("C101").Value = 0.000002145767211914062
$ws.Range("A102").Value = '''''que'''
$ws.Range("C102").Value = 0.000000476837158203125
$ws.Range("A103").Value = '''''que'''
$ws.Range("C103").Value = 0.00000095367431640625
$ws.Range("A104").Value = '''''que'''
$ws.Range("C104").Value = 0.000000476837158203125
$ws.Range("A105").Value = '''''que'''
$ws.Range("C105").Value = 0.0000007152557373046875
$ws.Range("A106").Value = '''''que'''
$ws.Range("C106").Value = 0.0000007152557373046875
$ws.Range("A107").Value = '''''que'''
$ws.Range("C107").Value = 0.000000476837158203125
$ws.Range("A108").Value = '''''que'''
$ws.Range("C108").Value = 0.000000476837158203125
$ws.Range("A109").Value = '''''que'''
$ws.Range("C109").Value = 0.000000476837158203125
$ws.Range("A110").Value = '''''que'''
$ws.Range("C110").Value = 0.000000476837158203125
$ws.Range("A111").Value = '''''que'''
$ws.Range("A112").Value = '''''que'''
$ws.Range("C112").Value = 0.000001192092895507812
$ws.Range("A113").Value = '''''que'''
